# "Updated points + bugfix"
#
# Sprint 3 table (rows 34-47) tracks daily "Points completed" in column B.
# Column F is a running cumulative total (F[r] = B[r] + F[r-1]) and
# column G is the points still left (G[r] = $E$34 - F[r]); both are driven
# by shared formulas already present in the sheet, so changing B43 alone
# makes Excel recalculate F43:F47 and G43:G47 automatically.
#
# Day 43 (2019-12-09) had 0 points logged; the bugfix corrects it to 5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B43").Value = 5

# Leave the sheet with B44 selected, matching the author's last cursor
# position when the workbook was saved.
$ws.Range("B44").Select()
